# Error.xlsx edit: "make SceneWarp relay on level to work"
#
# 1. Row Id=2000 ("不符合地图进入条件" / does-not-meet-map-entry-condition) gets
#    its description text changed to a level-gated message ("需要达到等级{0}才能进入").
# 2. A brand new row is inserted (Id=2002, "请先探索地图激活传送门") right before the
#    trailing "钻石不足" (Id=3000) row, which shifts down by one row.
# 3. The backing table / dimension / selection all grow to match the extra row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Re-word the existing level-gate error message (row 17, Id 2000).
$ws.Range("B17").Value = "需要达到等级{0}才能进入"

# 2. Insert the new row's data in what is currently the last row (19), then
#    re-create the old last row (Id 3000 / "钻石不足") one row below it.
$ws.Range("A19").Value = 2002
$ws.Range("B19").Value = "请先探索地图激活传送门"
$ws.Range("A20").Value = 3000
$ws.Range("B20").Value = "钻石不足"

# 3. Grow the table ("表1") so the new row is included.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B20"))

# 4. Match the author's final selection (their cursor ended up on the new
#    last cell after typing the extra row).
$ws.Range("B20").Select() | Out-Null
